$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# --- 1. Mark the header row of the Solutions table as a repeating table header ---
# (<w:trPr><w:tblHeader w:val="true"/></w:trPr> on the first row)
$tbl.Rows(1).HeadingFormat = $true

# --- 2. Normalize the table's preferred width to a clean 100% (pct) value ---
# (<w:tblW w:type="pct" w:w="5000"/> instead of the float "5000.0")
$tbl.PreferredWidthType = 2   # wdPreferredWidthPercent
$tbl.PreferredWidth = 250

# --- 3. Update the Oreo-filling confidence-interval answer (problem 13) ---
# Find the row whose "Problem" column reads "13"
$targetRow = $null
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
  if ($tbl.Rows($i).Cells(1).Range.Text -like "13*") {
    $targetRow = $i
    break
  }
}

if ($targetRow -ne $null) {
  $solutionCell = $tbl.Rows($targetRow).Cells(3)

  # Replace the bounds within that cell only, so the change stays local
  # to problem 13's solution cell.
  $solutionCell.Range.Find.Execute("2.808", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "2.535", 2)
  $solutionCell.Range.Find.Execute("2.988", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "3.165", 2)
}
